$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 135
$ws.Range("H2").Value = 'bedrooms'
$ws.Range("I2").Value = 'distractor'
$ws.Range("K2").Value = 'f'
$ws.Range("L2").Value = 'stimuli/img_fea1z.png'
$ws.Range("M2").Value = 79.45945945945945
$ws.Range("N2").Value = 56.24324324324324
$ws.Range("O2").Value = 67.85135135135135
$ws.Range("P2").Value = 37
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 7
$ws.Range("F3").Value = 136
$ws.Range("H3").Value = 'living_rooms'
$ws.Range("I3").Value = 'target'
$ws.Range("K3").Value = 'j'
$ws.Range("L3").Value = 'stimuli/img_xu1p3.png'
$ws.Range("M3").Value = 75.27659574468085
$ws.Range("N3").Value = 56.68085106382978
$ws.Range("O3").Value = 65.97872340425532
$ws.Range("P3").Value = 47
$ws.Range("Q3").Value = 7
$ws.Range("R3").Value = 7
$ws.Range("S3").Value = 7
$ws.Range("F4").Value = 137
$ws.Range("L4").Value = 'stimuli/img_4o8l0.png'
$ws.Range("M4").Value = 46.02173913043478
$ws.Range("N4").Value = 31.45652173913043
$ws.Range("O4").Value = 38.73913043478261
$ws.Range("P4").Value = 46
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = 3
$ws.Range("F5").Value = 138
$ws.Range("H5").Value = 'kitchens'
$ws.Range("I5").Value = 'distractor'
$ws.Range("K5").Value = 'f'
$ws.Range("L5").Value = 'stimuli/img_cxpff.png'
$ws.Range("M5").Value = 74.92307692307692
$ws.Range("N5").Value = 53.28205128205128
$ws.Range("O5").Value = 64.1025641025641
$ws.Range("P5").Value = 39
$ws.Range("F6").Value = 139
$ws.Range("L6").Value = 'stimuli/img_pey7u.png'
$ws.Range("M6").Value = 30.34883720930232
$ws.Range("N6").Value = 20.34883720930232
$ws.Range("O6").Value = 25.34883720930232
$ws.Range("P6").Value = 43
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 2
$ws.Range("F7").Value = 140
$ws.Range("L7").Value = 'stimuli/img_bbs77.png'
$ws.Range("M7").Value = 31.64444444444445
$ws.Range("N7").Value = 21.26666666666667
$ws.Range("O7").Value = 26.45555555555556
$ws.Range("P7").Value = 45
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 2
$ws.Range("F8").Value = 141
$ws.Range("L8").Value = 'stimuli/img_di6f0.png'
$ws.Range("M8").Value = 94.04347826086956
$ws.Range("N8").Value = 83.34782608695652
$ws.Range("O8").Value = 88.69565217391303
$ws.Range("P8").Value = 46
$ws.Range("Q8").Value = 10
$ws.Range("R8").Value = 10
$ws.Range("S8").Value = 10
$ws.Range("F9").Value = 142
$ws.Range("L9").Value = 'stimuli/img_0kqc0.png'
$ws.Range("M9").Value = 43.74468085106383
$ws.Range("N9").Value = 27.14893617021277
$ws.Range("O9").Value = 35.4468085106383
$ws.Range("P9").Value = 47
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 2
$ws.Range("F10").Value = 143
$ws.Range("L10").Value = 'stimuli/img_bj99b.png'
$ws.Range("M10").Value = 82.79069767441861
$ws.Range("N10").Value = 65.46511627906976
$ws.Range("O10").Value = 74.12790697674419
$ws.Range("Q10").Value = 8
$ws.Range("R10").Value = 8
$ws.Range("S10").Value = 8
$ws.Range("F11").Value = 144
$ws.Range("L11").Value = 'stimuli/img_cehin.png'
$ws.Range("M11").Value = 78.86363636363636
$ws.Range("N11").Value = 60.02272727272727
$ws.Range("O11").Value = 69.44318181818181
$ws.Range("P11").Value = 44
$ws.Range("Q11").Value = 7
$ws.Range("R11").Value = 7
$ws.Range("S11").Value = 7
$ws.Range("F12").Value = 145
$ws.Range("H12").Value = 'living_rooms'
$ws.Range("I12").Value = 'target'
$ws.Range("K12").Value = 'j'
$ws.Range("L12").Value = 'stimuli/img_amsgw.png'
$ws.Range("M12").Value = 86.08510638297872
$ws.Range("N12").Value = 65.95744680851064
$ws.Range("O12").Value = 76.02127659574468
$ws.Range("P12").Value = 47
$ws.Range("Q12").Value = 9
$ws.Range("R12").Value = 9
$ws.Range("S12").Value = 9
$ws.Range("F13").Value = 146
$ws.Range("H13").Value = 'bedrooms'
$ws.Range("I13").Value = 'distractor'
$ws.Range("K13").Value = 'f'
$ws.Range("L13").Value = 'stimuli/img_ys3qz.png'
$ws.Range("M13").Value = 46.79545454545455
$ws.Range("N13").Value = 31.20454545454545
$ws.Range("O13").Value = 39
$ws.Range("P13").Value = 44
$ws.Range("Q13").Value = 2
$ws.Range("F14").Value = 147
$ws.Range("H14").Value = 'living_rooms'
$ws.Range("I14").Value = 'target'
$ws.Range("K14").Value = 'j'
$ws.Range("L14").Value = 'stimuli/img_kost0.png'
$ws.Range("M14").Value = 63.09090909090909
$ws.Range("N14").Value = 42.77272727272727
$ws.Range("O14").Value = 52.93181818181819
$ws.Range("P14").Value = 44
$ws.Range("Q14").Value = 5
$ws.Range("R14").Value = 5
$ws.Range("S14").Value = 5
$ws.Range("F15").Value = 148
$ws.Range("L15").Value = 'stimuli/img_6a0hu.png'
$ws.Range("M15").Value = 61.275
$ws.Range("N15").Value = 42.025
$ws.Range("O15").Value = 51.65
$ws.Range("P15").Value = 40
$ws.Range("Q15").Value = 4
$ws.Range("R15").Value = 4
$ws.Range("S15").Value = 4
$ws.Range("F16").Value = 149
$ws.Range("H16").Value = 'kitchens'
$ws.Range("L16").Value = 'stimuli/img_kugyw.png'
$ws.Range("M16").Value = 74.25
$ws.Range("N16").Value = 54.10714285714285
$ws.Range("O16").Value = 64.17857142857143
$ws.Range("P16").Value = 28
$ws.Range("F17").Value = 150
$ws.Range("H17").Value = 'living_rooms'
$ws.Range("I17").Value = 'target'
$ws.Range("K17").Value = 'j'
$ws.Range("L17").Value = 'stimuli/img_xy930.png'
$ws.Range("M17").Value = 70.5952380952381
$ws.Range("N17").Value = 49.47619047619047
$ws.Range("O17").Value = 60.03571428571429
$ws.Range("F18").Value = 151
$ws.Range("L18").Value = 'stimuli/img_w8yhd.png'
$ws.Range("M18").Value = 55.74418604651163
$ws.Range("N18").Value = 38.90697674418605
$ws.Range("O18").Value = 47.32558139534883
$ws.Range("P18").Value = 43
$ws.Range("Q18").Value = 4
$ws.Range("R18").Value = 4
$ws.Range("S18").Value = 4
$ws.Range("F19").Value = 152
$ws.Range("L19").Value = 'stimuli/img_wgkqa.png'
$ws.Range("M19").Value = 87.25581395348837
$ws.Range("N19").Value = 71.13953488372093
$ws.Range("O19").Value = 79.19767441860465
$ws.Range("P19").Value = 43
$ws.Range("Q19").Value = 10
$ws.Range("R19").Value = 10
$ws.Range("S19").Value = 10
$ws.Range("F20").Value = 153
$ws.Range("L20").Value = 'stimuli/img_abobq.png'
$ws.Range("M20").Value = 75.1842105263158
$ws.Range("N20").Value = 54.13157894736842
$ws.Range("O20").Value = 64.65789473684211
$ws.Range("P20").Value = 38
$ws.Range("Q20").Value = 6
$ws.Range("R20").Value = 6
$ws.Range("S20").Value = 6
$ws.Range("F21").Value = 154
$ws.Range("L21").Value = 'stimuli/img_6zz63.png'
$ws.Range("M21").Value = 87.66666666666667
$ws.Range("N21").Value = 70.6
$ws.Range("O21").Value = 79.13333333333333
$ws.Range("P21").Value = 45
$ws.Range("Q21").Value = 9
$ws.Range("R21").Value = 10
$ws.Range("S21").Value = 10
$ws.Range("F22").Value = 155
$ws.Range("L22").Value = 'stimuli/img_eh0no.png'
$ws.Range("M22").Value = 53.66666666666666
$ws.Range("N22").Value = 36.02564102564103
$ws.Range("O22").Value = 44.84615384615385
$ws.Range("P22").Value = 39
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 3
$ws.Range("S22").Value = 3
$ws.Range("F23").Value = 156
$ws.Range("H23").Value = 'kitchens'
$ws.Range("L23").Value = 'stimuli/img_pt3d7.png'
$ws.Range("M23").Value = 65.08571428571429
$ws.Range("N23").Value = 44.65714285714286
$ws.Range("O23").Value = 54.87142857142857
$ws.Range("P23").Value = 35
$ws.Range("Q23").Value = 4
$ws.Range("R23").Value = 4
$ws.Range("S23").Value = 4
$ws.Range("F24").Value = 157
$ws.Range("L24").Value = 'stimuli/img_16kib.png'
$ws.Range("M24").Value = 80.97727272727273
$ws.Range("N24").Value = 61.11363636363637
$ws.Range("O24").Value = 71.04545454545455
$ws.Range("Q24").Value = 8
$ws.Range("R24").Value = 8
$ws.Range("S24").Value = 8
$ws.Range("F25").Value = 158
$ws.Range("L25").Value = 'stimuli/img_wz6x5.png'
$ws.Range("M25").Value = 68.3695652173913
$ws.Range("N25").Value = 48.47826086956522
$ws.Range("O25").Value = 58.42391304347826
$ws.Range("P25").Value = 46
$ws.Range("Q25").Value = 5
$ws.Range("R25").Value = 5
$ws.Range("S25").Value = 5
$ws.Range("F26").Value = 159
$ws.Range("F27").Value = 160
$ws.Range("H27").Value = 'bedrooms'
$ws.Range("I27").Value = 'distractor'
$ws.Range("K27").Value = 'f'
$ws.Range("L27").Value = 'stimuli/img_twj5p.png'
$ws.Range("M27").Value = 67.71739130434783
$ws.Range("N27").Value = 42.08695652173913
$ws.Range("O27").Value = 54.90217391304348
$ws.Range("P27").Value = 46
$ws.Range("Q27").Value = 4
$ws.Range("R27").Value = 4
$ws.Range("S27").Value = 4
